$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.420.03"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.725.97"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.32"
$ws.Range("E5").Value = "  -0.98%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("E7").Value = "  +1.45%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2621"
$ws.Range("E8").Value = "  -2.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06195"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "1.729.82"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07023"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.569"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5994"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.41"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "26.434.95"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007124"
$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -2.16%  "

$ws.Range("D21").Value = "1.953.47"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.492"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.601"
$ws.Range("E23").Value = "  -3.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.189"
$ws.Range("E24").Value = "  -1.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.01"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("E26").Value = "  -0.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.405"
$ws.Range("E27").Value = "  +0.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "107.02"
$ws.Range("E28").Value = "  -0.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.726"
$ws.Range("E29").Value = "  -3.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.968"
$ws.Range("E30").Value = "  -0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07984"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.688"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04526"
$ws.Range("E33").Value = "  -1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"

$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6237"
$ws.Range("E36").Value = "  -2.15%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9092"
$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.982"
$ws.Range("E38").Value = "  -5.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.406"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01491"
$ws.Range("E41").Value = "  -0.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.25"
$ws.Range("E42").Value = "  -4.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.401"
$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3870"
$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.707"
$ws.Range("E45").Value = "  -3.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1157"
$ws.Range("E46").Value = "  -2.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05365"
$ws.Range("E47").Value = "  +0.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.30"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.703"
$ws.Range("E49").Value = "  -2.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.253"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.11"
$ws.Range("E51").Value = "  -0.45%  "
